$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the source value and let the dependent formula (F4 = E4/60) recalc
$ws.Range("E4").Value = 5124

# Update the active selection to E4 (matches the saved cursor position in the diff)
$ws.Range("E4").Select()

# Ensure the workbook recalculates so F4 reflects the new value
$excel.Calculate()
